$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Cod. Moneda" column (E) for the detail rows was recorded as "US$";
# update it to "USD" across all the data rows (2-19).
$ws.Range("E2:E19").Value = "USD"
